# Remove column from alcohol data
# Column M (the penultimate data column) is deleted entirely; the data
# that used to live in column N shifts left to become the new column M.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the entire column M - this shifts column N (and anything to its
# right) one column to the left, so the old N values become the new M
# values and the old M values are gone.
$ws.Range("M1").EntireColumn.Delete() | Out-Null

# Move the active selection to M4 (the new rightmost data column at the
# previously visible row), matching the post-edit selection.
$ws.Range("M4").Select() | Out-Null
